$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "savedEachDay" (sheet2.xml) - add rows 100-110
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("savedEachDay")
$ws.Activate()

# Row 100 - pTranspirationEfficiencyLinkedToCO2
$ws.Cells.Item(100, 1).Value = "pTranspirationEfficiencyLinkedToCO2"
$ws.Cells.Item(100, 2).Value = "CropParameter"
$ws.Cells.Item(100, 3).Value = "numeric"
$ws.Cells.Item(100, 4).Value = "rWaterBudget"
$ws.Cells.Item(100, 5).Value = "-"
$ws.Cells.Item(100, 6).Value = "Transpiration efficiency coefficient at CO2 concentration 350 ppm"
$ws.Cells.Item(100, 7).Value = "TEC350"
$ws.Cells.Item(100, 8).Value = "coefficient d'efficacité de la transpiration lié à la teneur en CO2"
$ws.Cells.Item(100, 9).Formula = "=NA()"

# Row 101 - cDownwardWaterFlux.1
$ws.Cells.Item(101, 1).Value = "cDownwardWaterFlux.1"
$ws.Cells.Item(101, 2).Value = "computed"
$ws.Cells.Item(101, 3).Value = "numeric"
$ws.Cells.Item(101, 4).Value = "rWaterBudget"
$ws.Cells.Item(101, 5).Value = "mm"
$ws.Cells.Item(101, 7).Value = "FLOUT"

# Rows 102-110 - cDownwardWaterFlux.2 .. cDownwardWaterFlux.10
for ($i = 2; $i -le 10; $i++) {
    $r = 100 + $i
    $ws.Cells.Item($r, 1).Value = "cDownwardWaterFlux.$i"
    $ws.Cells.Item($r, 2).Value = "computed"
    $ws.Cells.Item($r, 3).Value = "numeric"
    $ws.Cells.Item($r, 4).Value = "rWaterBudget"
    $ws.Cells.Item($r, 5).Value = "mm"
}

# Update the active selection to match the post-edit state
$ws.Range("F101").Select()

# ---------------------------------------------------------------------------
# Sheet "other" (sheet3.xml) - add row 9
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("other")
$ws3.Activate()

$ws3.Cells.Item(9, 1).Value = "pVPDcoef"
$ws3.Cells.Item(9, 2).Value = "SoilParameter"
$ws3.Cells.Item(9, 3).Value = "numeric"
$ws3.Cells.Item(9, 4).Value = "rWaterBudget"
$ws3.Cells.Item(9, 6).Value = "A coefficient to calculate VPD; 0.65 for humid and subhumid climates and 0.75 for arid and semi-arid climates"
$ws3.Cells.Item(9, 7).Value = "VPDF"
$ws3.Cells.Item(9, 8).Value = "coefficient pour calculer le VPD, 0.65 pour les climats humides et subhumides, 0.75 pour les climats arides et semiarides"

$ws3.Range("D37").Select()

# Re-activate the sheet that was originally selected/visible
$ws.Activate()
